# Refresh the cryptocurrency price/volume figures with the latest scrape,
# and fix the ranking order so row 40 is Kaspa and row 41 is ARBITRUM
# (each keeping its own coin's link/price/volume data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose Price ("D") is numeric-looking are written with a leading
# apostrophe so Excel keeps them as text (matching the source data, which
# stores prices as plain strings) instead of auto-converting to a number.
$updates = @(
    @{Row=2; Col="D"; Value="43.064.03"},
    @{Row=2; Col="E"; Value="  +1.06%  "},
    @{Row=3; Col="D"; Value="2.302.90"},
    @{Row=3; Col="E"; Value="  +0.81%  "},
    @{Row=4; Col="E"; Value="  -0.03%  "},
    @{Row=5; Col="D"; Value="'300.18"},
    @{Row=5; Col="E"; Value="  -0.15%  "},
    @{Row=6; Col="D"; Value="'97.69"},
    @{Row=6; Col="E"; Value="  +0.62%  "},
    @{Row=7; Col="E"; Value="  +0.55%  "},
    @{Row=8; Col="E"; Value="  -0.02%  "},
    @{Row=9; Col="D"; Value="'0.506"},
    @{Row=9; Col="E"; Value="  +1.60%  "},
    @{Row=10; Col="D"; Value="'33.82"},
    @{Row=10; Col="E"; Value="  +1.94%  "},
    @{Row=11; Col="E"; Value="  +0.63%  "},
    @{Row=12; Col="D"; Value="'49.08"},
    @{Row=12; Col="E"; Value="  -0.42%  "},
    @{Row=13; Col="D"; Value="'0.116"},
    @{Row=13; Col="E"; Value="  +2.98%  "},
    @{Row=14; Col="D"; Value="'17.19"},
    @{Row=14; Col="E"; Value="  +11.73%  "},
    @{Row=15; Col="D"; Value="'6.79"},
    @{Row=15; Col="E"; Value="  +2.07%  "},
    @{Row=16; Col="D"; Value="2.660.76"},
    @{Row=16; Col="E"; Value="  +0.78%  "},
    @{Row=17; Col="D"; Value="2.295.13"},
    @{Row=17; Col="E"; Value="  -1.28%  "},
    @{Row=18; Col="D"; Value="'0.811"},
    @{Row=18; Col="E"; Value="  +3.13%  "},
    @{Row=19; Col="D"; Value="43.002.49"},
    @{Row=19; Col="E"; Value="  +1.09%  "},
    @{Row=20; Col="D"; Value="'11.67"},
    @{Row=20; Col="E"; Value="  +1.76%  "},
    @{Row=21; Col="E"; Value="  +0.87%  "},
    @{Row=22; Col="D"; Value="'6.05"},
    @{Row=22; Col="E"; Value="  +0.71%  "},
    @{Row=23; Col="E"; Value="  +1.08%  "},
    @{Row=24; Col="D"; Value="'236.57"},
    @{Row=25; Col="E"; Value="  +5.44%  "},
    @{Row=27; Col="E"; Value="  -1.52%  "},
    @{Row=28; Col="D"; Value="'24.39"},
    @{Row=28; Col="E"; Value="  +0.13%  "},
    @{Row=29; Col="D"; Value="'166.52"},
    @{Row=29; Col="E"; Value="  +0.83%  "},
    @{Row=30; Col="D"; Value="'2.08"},
    @{Row=30; Col="E"; Value="  +0.45%  "},
    @{Row=31; Col="D"; Value="'33.82"},
    @{Row=31; Col="E"; Value="  -0.26%  "},
    @{Row=32; Col="E"; Value="  +0.31%  "},
    @{Row=33; Col="E"; Value="  +0.01%  "},
    @{Row=34; Col="E"; Value="  +0.22%  "},
    @{Row=35; Col="E"; Value="  +6.36%  "},
    @{Row=36; Col="E"; Value="  +1.93%  "},
    @{Row=37; Col="D"; Value="'16.91"},
    @{Row=37; Col="E"; Value="  +4.31%  "},
    @{Row=38; Col="D"; Value="'0.0703"},
    @{Row=38; Col="E"; Value="  +1.04%  "},
    @{Row=39; Col="E"; Value="  +0.09%  "},
    @{Row=42; Col="E"; Value="  -0.09%  "},
    @{Row=44; Col="D"; Value="1.995.81"},
    @{Row=44; Col="E"; Value="  +1.74%  "},
    @{Row=45; Col="D"; Value="'0.0284"},
    @{Row=45; Col="E"; Value="  +0.74%  "},
    @{Row=46; Col="D"; Value="'9.86"},
    @{Row=46; Col="E"; Value="  +1.88%  "},
    @{Row=47; Col="D"; Value="'17.50"},
    @{Row=47; Col="E"; Value="  -2.30%  "},
    @{Row=48; Col="E"; Value="  +1.14%  "},
    @{Row=49; Col="D"; Value="2.527.39"},
    @{Row=49; Col="E"; Value="  +0.72%  "},
    @{Row=50; Col="D"; Value="'53.39"},
    @{Row=50; Col="E"; Value="  +0.40%  "},
    @{Row=51; Col="E"; Value="  -1.59%  "}
)

foreach ($u in $updates) {
    $cellRef = "$($u.Col)$($u.Row)"
    $ws.Range($cellRef).Value = $u.Value
}

# Rows 40 and 41 swap coins (ARBITRUM <-> Kaspa) along with their data.
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.101"
$ws.Range("E40").Value = "  +1.95%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "  +1.24%  "
